# edit.ps1 - apply the "add some examples on cassandra" commit
#
# 1) Bump the cached fixed date/time fields from 20/12/2021 (or 12/20/2021)
#    to 21/12/2021 (or 12/21/2021) across the slide master, all slide
#    layouts, and the notes master.
# 2) Tighten the "PRIMARY KEY ( ... ) );" Cassandra CQL snippets on slide 25
#    to "PRIMARY KEY ( ... ));" (drop the stray space before the final
#    closing parens).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date / time placeholder fields
# ---------------------------------------------------------------------

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = ""
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ("$phType" -eq "16") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster

# Slide master itself (MM/DD/YYYY style cached field)
Set-DatePlaceholderText $master.Shapes "12/21/2021"

# Every slide layout off this master (also MM/DD/YYYY style)
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "12/21/2021"
}

# Notes master uses the Italian DD/MM/YYYY cached field and only responds
# to the header/footer Date-and-Time object (direct shape text edits are
# ignored there).
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "21/12/2021"

# ---------------------------------------------------------------------
# 2) Slide 25 ("Morale") Cassandra CQL code boxes
# ---------------------------------------------------------------------

$slide = $p.Slides.Item(25)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text

    if ($full.IndexOf(" ((year), id) );") -ge 0) {
        $start = $full.IndexOf(" ((year), id) );") + 1
        $len = " ((year), id) );".Length
        $sub = $tr.Characters($start, $len)
        $sub.Text = " ((year), id));"
    }
    elseif ($full.IndexOf(" ((year), name, id) );") -ge 0) {
        $start = $full.IndexOf(" ((year), name, id) );") + 1
        $len = " ((year), name, id) );".Length
        $sub = $tr.Characters($start, $len)
        $sub.Text = " ((year), name, id));"
    }
    elseif ($full.IndexOf(" (id) );") -ge 0) {
        $start = $full.IndexOf(" (id) );") + 1
        $len = " (id) );".Length
        $sub = $tr.Characters($start, $len)
        $sub.Text = " (id));"
    }
}
